$wb = $excel.ActiveWorkbook

$wsUni = $wb.Worksheets.Item("Университеты")
$wsUni.Range("A3").Value = "0020-high"
$wsUni.Range("E5").Value = "MEDICINE"

$wsStud = $wb.Worksheets.Item("Студенты")
$wsStud.Activate()
$wsStud.Range("B2").Select()
